# Horarios Linea 141 - scrape refresh (05:18:23 -> 05:47:32)
# Updates the three timetable sheets (LP1912, LP1912-215, 6203-6173) with the
# newly scraped rows: two new "just departed" rows inserted near the top of
# the still-pending LP1912 entries, one additional LP1912 row inserted
# further down, seven brand new LP1912 rows appended at the end, and one
# brand new row appended to 6203-6173.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "LP1912"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 05:47:32"
$ws1.Range("A3").Value = "Total filas: 49"

# Final state (rows 28-54) after the two inserted rows near the top, the
# one inserted row further down, and the seven appended rows at the end.
$rows1 = @(
    @("05:47:32","05:47","15_ABASTO",0,"LP1912"),
    @("05:47:32","05:49","14_ABASTO",2,"LP1912"),
    @("04:40:32","06:04","16_SANTA ANA",84,"LP1912"),
    @("04:18:06","06:09","16_SANTA ANA",111,"LP1912"),
    @("04:40:32","06:11","215A_EL PATO",91,"LP1912"),
    @("04:18:06","06:12","215A_EL PATO",114,"LP1912"),
    @("04:18:06","06:14","225_HARAS DEL SUR",116,"LP1912"),
    @("04:40:32","06:21","26_HERNANDEZ",101,"LP1912"),
    @("04:40:32","06:27","23_HERNANDEZ",107,"LP1912"),
    @("04:40:32","06:29","86_EST CHICA-ESC AGRARIA",109,"LP1912"),
    @("04:40:32","06:31","16_SANTA ANA",111,"LP1912"),
    @("04:53:50","06:44","225_C ROCA-H SUR",111,"LP1912"),
    @("04:53:50","06:46","215C_EL PATO",113,"LP1912"),
    @("05:18:23","06:58","10_OLMOS",100,"LP1912"),
    @("05:18:23","06:59","14_ABASTO",101,"LP1912"),
    @("05:47:32","07:04","23_HERNANDEZ",77,"LP1912"),
    @("05:18:23","07:05","15_ABASTO",107,"LP1912"),
    @("05:18:23","07:07","225_GOMEZ",109,"LP1912"),
    @("05:18:23","07:11","215A_EL PATO",113,"LP1912"),
    @("05:18:23","07:15","11_ETCHEVERRY",117,"LP1912"),
    @("05:47:32","07:21","26_HERNANDEZ",94,"LP1912"),
    @("05:47:32","07:27","10_OLMOS",100,"LP1912"),
    @("05:47:32","07:31","11_ETCHEVERRY",104,"LP1912"),
    @("05:47:32","07:31","16_SANTA ANA",104,"LP1912"),
    @("05:47:32","07:32","84_COLONIA URQUIZA-ESC 49",105,"LP1912"),
    @("05:47:32","07:36","27_EL RETIRO",109,"LP1912"),
    @("05:47:32","07:39","10_OLMOS",112,"LP1912")
)

$startRow1 = 28
for ($i = 0; $i -lt $rows1.Length; $i++) {
    $r = $startRow1 + $i
    $data = $rows1[$i]
    $ws1.Cells.Item($r, 1).Value = $data[0]
    $ws1.Cells.Item($r, 2).Value = $data[1]
    $ws1.Cells.Item($r, 3).Value = $data[2]
    $ws1.Cells.Item($r, 4).Value = $data[3]
    $ws1.Cells.Item($r, 5).Value = $data[4]
}

# ---------------------------------------------------------------------
# Sheet "LP1912-215"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 05:47:32"

# ---------------------------------------------------------------------
# Sheet "6203-6173"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 05:47:32"
$ws3.Range("A3").Value = "Total filas: 9"

$ws3.Cells.Item(14, 1).Value = "05:47:32"
$ws3.Cells.Item(14, 2).Value = "07:35"
$ws3.Cells.Item(14, 3).Value = "215A_LA PLATA"
$ws3.Cells.Item(14, 4).Value = 108
$ws3.Cells.Item(14, 5).Value = "L6173"
